$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is protected; temporarily unprotect so the cells can be edited,
# then restore protection (same effective settings: contents/objects/
# scenarios protected, formatting of columns/rows not allowed) afterwards.
$ws.Unprotect()

# --- Update the "as of" date in the confidential disclosure note (A11) ---
[void]$ws.Cells.Replace("2021-04-21", "2021-04-22")

# --- Refresh the model holdings weight / percent-change figures ---
$ws.Range("D2").Value = 0.4918669639687679
$ws.Range("E2").Value = -0.007766856828648083

$ws.Range("D3").Value = 0.2513442368008163
$ws.Range("E3").Value = -0.01034482758620692

$ws.Range("D4").Value = 0.09883367576347579
$ws.Range("E4").Value = -0.002448879637565904

$ws.Range("D5").Value = 0.1013652128880768
$ws.Range("E5").Value = -0.00606405154443812

$ws.Range("D6").Value = 0.02916771034866406
$ws.Range("E6").Value = -0.006519806381507509

$ws.Range("D7").Value = 0.02742220023019911
$ws.Range("E7").Value = -0.003313299429804295

$ws.Range("D8").Value = 1
$ws.Range("E8").Value = -0.007558114518679382

# Restore sheet protection (contents/objects/scenarios locked, matching the
# original protection flags).
$ws.Protect($null, $true, $true, $true)
